$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (row 1, columns B..R): insert "企业" before "电子商务销售额"
$headerCols = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)  # B..R
foreach ($c in $headerCols) {
    $cell = $ws.Cells.Item(1, $c)
    $oldText = $cell.Value()
    $newText = $oldText.Replace("电子商务销售额", "企业电子商务销售额")
    $cell.Value = $newText
}

# 2. Move column O (the grand-total "企业电子商务销售额" column, now holding the updated
#    header text) to become column C, shifting old C..N right by one into D..O.
#    Done cell-by-cell (row by row) across the used rows (1..9) so no column-level
#    metadata (widths etc.) gets stamped onto the sheet.
for ($r = 1; $r -le 9; $r++) {
    $oldC = $ws.Cells.Item($r, 3).Value()
    $oldD = $ws.Cells.Item($r, 4).Value()
    $oldE = $ws.Cells.Item($r, 5).Value()
    $oldF = $ws.Cells.Item($r, 6).Value()
    $oldG = $ws.Cells.Item($r, 7).Value()
    $oldH = $ws.Cells.Item($r, 8).Value()
    $oldI = $ws.Cells.Item($r, 9).Value()
    $oldJ = $ws.Cells.Item($r, 10).Value()
    $oldK = $ws.Cells.Item($r, 11).Value()
    $oldL = $ws.Cells.Item($r, 12).Value()
    $oldM = $ws.Cells.Item($r, 13).Value()
    $oldN = $ws.Cells.Item($r, 14).Value()
    $oldO = $ws.Cells.Item($r, 15).Value()

    $ws.Cells.Item($r, 3).Value = $oldO
    $ws.Cells.Item($r, 4).Value = $oldC
    $ws.Cells.Item($r, 5).Value = $oldD
    $ws.Cells.Item($r, 6).Value = $oldE
    $ws.Cells.Item($r, 7).Value = $oldF
    $ws.Cells.Item($r, 8).Value = $oldG
    $ws.Cells.Item($r, 9).Value = $oldH
    $ws.Cells.Item($r, 10).Value = $oldI
    $ws.Cells.Item($r, 11).Value = $oldJ
    $ws.Cells.Item($r, 12).Value = $oldK
    $ws.Cells.Item($r, 13).Value = $oldL
    $ws.Cells.Item($r, 14).Value = $oldM
    $ws.Cells.Item($r, 15).Value = $oldN
}

# 3. Append new row 10 for "2021年"
$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 9421.911599999999
$ws.Range("C10").Value = 227611.27949
$ws.Range("D10").Value = 1894.5493
$ws.Range("E10").Value = 18604.47603
$ws.Range("F10").Value = 71267.6548
$ws.Range("G10").Value = 45.22657
$ws.Range("H10").Value = 122.05253
$ws.Range("I10").Value = 208.72657
$ws.Range("J10").Value = 370.56701
$ws.Range("K10").Value = 115344.04555
$ws.Range("L10").Value = 800.38899
$ws.Range("M10").Value = 482.53263
$ws.Range("N10").Value = 72.8618
$ws.Range("O10").Value = 1799.74623
$ws.Range("P10").Value = 665.6768499999999
$ws.Range("Q10").Value = 5667.52455
$ws.Range("R10").Value = 843.33849

# Match the style of A2:A9 (style index 1, bold/bordered/centered) for the new A10 label cell
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
